# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $text) {
    # Force the cell to keep a literal text value even when it looks numeric
    # (e.g. '236.86', '44.40', '0.160'), matching the source data which stores
    # every Price/Volume column cell as a string, never a number.
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = 'Normal'
}

$ws.Cells.Item(2, 4).Value = '96.463.32'
$ws.Cells.Item(2, 5).Value = '  -0.25%  '

$ws.Cells.Item(3, 4).Value = '3.702.77'
$ws.Cells.Item(3, 5).Value = '  +1.54%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

Set-TextValue $ws.Cells.Item(5, 4) '236.86'
$ws.Cells.Item(5, 5).Value = '  -2.27%  '

Set-TextValue $ws.Cells.Item(6, 4) '1.89'
$ws.Cells.Item(6, 5).Value = '  +1.93%  '

Set-TextValue $ws.Cells.Item(7, 4) '652.29'
$ws.Cells.Item(7, 5).Value = '  -0.46%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.428'
$ws.Cells.Item(8, 5).Value = '  +1.46%  '

Set-TextValue $ws.Cells.Item(9, 4) '0.999'
$ws.Cells.Item(9, 5).Value = '  -0.02%  '

$ws.Cells.Item(10, 5).Value = '  -1.00%  '

$ws.Cells.Item(11, 4).Value = '3.699.68'
$ws.Cells.Item(11, 5).Value = '  +1.50%  '

Set-TextValue $ws.Cells.Item(12, 4) '44.40'
$ws.Cells.Item(12, 5).Value = '  +0.25%  '

$ws.Cells.Item(13, 5).Value = '  -0.05%  '

$ws.Cells.Item(14, 2).Value = 'ShibaInu'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(14, 4) '0.0000298'
$ws.Cells.Item(14, 5).Value = '  +15.38%  '

$ws.Cells.Item(15, 2).Value = 'Toncoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Cells.Item(15, 4) '6.78'
$ws.Cells.Item(15, 5).Value = '  +4.14%  '

$ws.Cells.Item(16, 4).Value = '4.389.62'
$ws.Cells.Item(16, 5).Value = '  +1.47%  '

$ws.Cells.Item(17, 4).Value = '96.205.15'
$ws.Cells.Item(17, 5).Value = '  -0.18%  '

Set-TextValue $ws.Cells.Item(18, 4) '8.84'
$ws.Cells.Item(18, 5).Value = '  +13.79%  '

$ws.Cells.Item(19, 4).Value = '3.698.78'
$ws.Cells.Item(19, 5).Value = '  +1.33%  '

Set-TextValue $ws.Cells.Item(20, 4) '13.03'
$ws.Cells.Item(20, 5).Value = '  +0.53%  '

Set-TextValue $ws.Cells.Item(21, 4) '18.90'
$ws.Cells.Item(21, 5).Value = '  +2.88%  '

Set-TextValue $ws.Cells.Item(22, 4) '0.505'
$ws.Cells.Item(22, 5).Value = '  -5.68%  '

Set-TextValue $ws.Cells.Item(23, 4) '518.22'
$ws.Cells.Item(23, 5).Value = '  +1.11%  '

Set-TextValue $ws.Cells.Item(24, 4) '3.38'
$ws.Cells.Item(24, 5).Value = '  -1.90%  '

Set-TextValue $ws.Cells.Item(25, 4) '0.0000206'
$ws.Cells.Item(25, 5).Value = '  +0.78%  '

$ws.Cells.Item(26, 5).Value = '  +1.25%  '

Set-TextValue $ws.Cells.Item(27, 4) '100.83'
$ws.Cells.Item(27, 5).Value = '  -0.32%  '

$ws.Cells.Item(28, 5).Value = '  +0.44%  '

$ws.Cells.Item(29, 5).Value = '  +3.22%  '

Set-TextValue $ws.Cells.Item(30, 4) '3.02'
$ws.Cells.Item(30, 5).Value = '  -0.55%  '

Set-TextValue $ws.Cells.Item(31, 4) '12.11'
$ws.Cells.Item(31, 5).Value = '  +1.93%  '

$ws.Cells.Item(32, 5).Value = '  +0.21%  '

$ws.Cells.Item(33, 5).Value = '  +5.58%  '

$ws.Cells.Item(34, 5).Value = '  -0.67%  '

Set-TextValue $ws.Cells.Item(35, 4) '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.23%  '

Set-TextValue $ws.Cells.Item(36, 4) '657.62'
$ws.Cells.Item(36, 5).Value = '  +6.86%  '

Set-TextValue $ws.Cells.Item(37, 4) '32.24'
$ws.Cells.Item(37, 5).Value = '  -2.72%  '

Set-TextValue $ws.Cells.Item(38, 4) '0.588'
$ws.Cells.Item(38, 5).Value = '  +0.60%  '

Set-TextValue $ws.Cells.Item(39, 4) '8.87'
$ws.Cells.Item(39, 5).Value = '  +0.55%  '

Set-TextValue $ws.Cells.Item(41, 4) '2.09'
$ws.Cells.Item(41, 5).Value = '  +7.51%  '

$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(42, 4) '41.39'
$ws.Cells.Item(42, 5).Value = '  -2.95%  '

$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(43, 4) '6.87'
$ws.Cells.Item(43, 5).Value = '  +11.65%  '

Set-TextValue $ws.Cells.Item(44, 4) '0.160'
$ws.Cells.Item(44, 5).Value = '  +0.84%  '

Set-TextValue $ws.Cells.Item(45, 4) '0.966'
$ws.Cells.Item(45, 5).Value = '  +1.46%  '

$ws.Cells.Item(46, 5).Value = '  +1.06%  '

Set-TextValue $ws.Cells.Item(47, 4) '0.434'
$ws.Cells.Item(47, 5).Value = '  +4.82%  '

Set-TextValue $ws.Cells.Item(48, 4) '2.28'
$ws.Cells.Item(48, 5).Value = '  -1.42%  '

Set-TextValue $ws.Cells.Item(49, 4) '23.57'

Set-TextValue $ws.Cells.Item(50, 4) '8.47'
$ws.Cells.Item(50, 5).Value = '  -1.70%  '

$ws.Cells.Item(51, 5).Value = '  +2.54%  '
